$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.050240391420051
$ws.Range("D2").Value = 1.049221199704704
$ws.Range("E2").Value = 1.063661699733451
$ws.Range("F2").Value = 1.070648850203086
$ws.Range("I2").Value = 1.039052624044855
$ws.Range("J2").Value = 1.055274591449345
$ws.Range("K2").Value = 1.051978862053855
$ws.Range("L2").Value = 1.066379734503542
$ws.Range("M2").Value = 1.073348125940836
$ws.Range("N2").Value = 1.021981237685588

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051873240476058
$ws.Range("D3").Value = 1.050468178356056
$ws.Range("E3").Value = 1.06532478135365
$ws.Range("F3").Value = 1.072459115922767
$ws.Range("I3").Value = 1.039464437775116
$ws.Range("J3").Value = 1.056553705720592
$ws.Range("K3").Value = 1.053036920329153
$ws.Range("L3").Value = 1.067855744475702
$ws.Range("M3").Value = 1.074972337630042
$ws.Range("N3").Value = 1.022429112652596

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052927248968394
$ws.Range("D4").Value = 1.051272542997792
$ws.Range("E4").Value = 1.066398720774397
$ws.Range("F4").Value = 1.073628473438989
$ws.Range("I4").Value = 1.039728086891468
$ws.Range("J4").Value = 1.057378488450905
$ws.Range("K4").Value = 1.053718488800514
$ws.Range("L4").Value = 1.0688081652814
$ws.Range("M4").Value = 1.076020861551118
$ws.Range("N4").Value = 1.022717374409033

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053369755151542
$ws.Range("D5").Value = 1.051610104539928
$ws.Range("E5").Value = 1.066849693748698
$ws.Range("F5").Value = 1.074119604472804
$ws.Range("I5").Value = 1.039838253060281
$ws.Range("J5").Value = 1.057724545123439
$ws.Range("K5").Value = 1.054004294156875
$ws.Range("L5").Value = 1.069207938175673
$ws.Range("M5").Value = 1.076461086702594
$ws.Range("N5").Value = 1.022838193669189

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053444019121868
$ws.Range("D6").Value = 1.051666748058735
$ws.Range("E6").Value = 1.066925384469973
$ws.Range("F6").Value = 1.074202040540039
$ws.Range("I6").Value = 1.039856711148353
$ws.Range("J6").Value = 1.05778260985318
$ws.Range("K6").Value = 1.054052239789354
$ws.Range("L6").Value = 1.069275025485805
$ws.Range("M6").Value = 1.076534969116141
$ws.Range("N6").Value = 1.022858458393989

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052933164094138
$ws.Range("D7").Value = 1.051277055831343
$ws.Range("E7").Value = 1.06640474868921
$ws.Range("F7").Value = 1.073635037772137
$ws.Range("I7").Value = 1.039729561572303
$ws.Range("J7").Value = 1.057383115143289
$ws.Range("K7").Value = 1.053722310587428
$ws.Range("L7").Value = 1.068813509506132
$ws.Range("M7").Value = 1.076026746101033
$ws.Range("N7").Value = 1.022718990235125

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050792756758147
$ws.Range("D8").Value = 1.049643147635481
$ws.Range("E8").Value = 1.064224205793682
$ws.Range("F8").Value = 1.071261061223271
$ws.Range("I8").Value = 1.039192384926026
$ws.Range("J8").Value = 1.055707478380969
$ws.Range("K8").Value = 1.052337077455547
$ws.Range("L8").Value = 1.066879115719281
$ws.Range("M8").Value = 1.073897550441697
$ws.Range("N8").Value = 1.022132920783579

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047000963699768
$ws.Range("D9").Value = 1.046744355888968
$ws.Range("E9").Value = 1.060364527432044
$ws.Range("F9").Value = 1.067061855352725
$ws.Range("I9").Value = 1.038224033477061
$ws.Range("J9").Value = 1.052732217650355
$ws.Range("K9").Value = 1.049872278000518
$ws.Range("L9").Value = 1.063449621096011
$ws.Range("M9").Value = 1.070126327042454
$ws.Range("N9").Value = 1.021088214646652

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044458784807293
$ws.Range("D10").Value = 1.044798085729369
$ws.Range("E10").Value = 1.057779023785287
$ws.Range("F10").Value = 1.064250819248978
$ws.Range("I10").Value = 1.03756359037473
$ws.Range("J10").Value = 1.050732909981039
$ws.Range("K10").Value = 1.048212541508934
$ws.Range("L10").Value = 1.061148571523736
$ws.Range("M10").Value = 1.067598427313242
$ws.Range("N10").Value = 1.020383470380639

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043354421140583
$ws.Range("D11").Value = 1.04395194849048
$ws.Range("E11").Value = 1.056656363511461
$ws.Range("F11").Value = 1.063030678969835
$ws.Range("I11").Value = 1.037274030246868
$ws.Range("J11").Value = 1.04986330046798
$ws.Range("K11").Value = 1.047489821032208
$ws.Range("L11").Value = 1.060148542536478
$ws.Range("M11").Value = 1.066500385219437
$ws.Range("N11").Value = 1.020076296423197

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04294365764126
$ws.Range("D12").Value = 1.043637135658261
$ws.Range("E12").Value = 1.056238873892908
$ws.Range("F12").Value = 1.062577006138884
$ws.Range("I12").Value = 1.037165931649251
$ws.Range("J12").Value = 1.049539691303719
$ws.Range("K12").Value = 1.047220753027432
$ws.Range("L12").Value = 1.059776523776032
$ws.Range("M12").Value = 1.066091991024479
$ws.Range("N12").Value = 1.019961891285515

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043031793107843
$ws.Range("D13").Value = 1.043704687771855
$ws.Range("E13").Value = 1.05632844895748
$ws.Range("F13").Value = 1.062674341507457
$ws.Range("I13").Value = 1.037189143815598
$ws.Range("J13").Value = 1.049609133800937
$ws.Range("K13").Value = 1.047278497129002
$ws.Range("L13").Value = 1.059856348782558
$ws.Range("M13").Value = 1.06617961730237
$ws.Range("N13").Value = 1.019986445540475

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043320478685172
$ws.Range("D14").Value = 1.043925936628637
$ws.Range("E14").Value = 1.056621863641211
$ws.Range("F14").Value = 1.062993187677972
$ws.Range("I14").Value = 1.037265105888611
$ws.Range("J14").Value = 1.049836563092388
$ws.Range("K14").Value = 1.047467592440471
$ws.Range("L14").Value = 1.060117802924023
$ws.Range("M14").Value = 1.066466638177041
$ws.Range("N14").Value = 1.020066845947744

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043498273571354
$ws.Range("D15").Value = 1.044062186156332
$ws.Range("E15").Value = 1.05680258158247
$ws.Range("F15").Value = 1.06318957809509
$ws.Range("I15").Value = 1.037311836526763
$ws.Range("J15").Value = 1.049976610249073
$ws.Range("K15").Value = 1.047584018186806
$ws.Range("L15").Value = 1.060278818372408
$ws.Range("M15").Value = 1.066643410204
$ws.Range("N15").Value = 1.020116342474819

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04453200098645
$ws.Range("D16").Value = 1.044854168761083
$ws.Range("E16").Value = 1.057853464049471
$ws.Range("F16").Value = 1.06433173255119
$ws.Range("I16").Value = 1.037582731623961
$ws.Range("J16").Value = 1.050790539982045
$ws.Range("K16").Value = 1.048260420047161
$ws.Range("L16").Value = 1.061214861890577
$ws.Range("M16").Value = 1.067671226901878
$ws.Range("N16").Value = 1.020403813682603

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045179460752697
$ws.Range("D17").Value = 1.045350043661484
$ws.Range("E17").Value = 1.058511809080899
$ws.Range("F17").Value = 1.065047376506471
$ws.Range("I17").Value = 1.037751694196334
$ws.Range("J17").Value = 1.051300044983873
$ws.Range("K17").Value = 1.048683619495685
$ws.Range("L17").Value = 1.061801027929631
$ws.Range("M17").Value = 1.068315016397158
$ws.Range("N17").Value = 1.020583594260085

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04555676837246
$ws.Range("D18").Value = 1.045638952792281
$ws.Range("E18").Value = 1.058895510562619
$ws.Range("F18").Value = 1.065464516657425
$ws.Range("I18").Value = 1.037849901727413
$ws.Range("J18").Value = 1.051596855540412
$ws.Range("K18").Value = 1.048930075100194
$ws.Range("L18").Value = 1.062142576522769
$ws.Range("M18").Value = 1.068690196557992
$ws.Range("N18").Value = 1.020688263101136

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045685362511988
$ws.Range("D19").Value = 1.045737408309739
$ws.Range("E19").Value = 1.059026292407227
$ws.Range("F19").Value = 1.065606703179576
$ws.Range("I19").Value = 1.037883329497535
$ws.Range("J19").Value = 1.05169799701277
$ws.Range("K19").Value = 1.049014044329405
$ws.Range("L19").Value = 1.062258976449485
$ws.Range("M19").Value = 1.068818067541979
$ws.Range("N19").Value = 1.020723919673386

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045110030212463
$ws.Range("D20").Value = 1.045296874805593
$ws.Range("E20").Value = 1.058441206030962
$ws.Range("F20").Value = 1.064970624057721
$ws.Range("I20").Value = 1.03773360187931
$ws.Range("J20").Value = 1.05124541880409
$ws.Range("K20").Value = 1.048638254556824
$ws.Range("L20").Value = 1.061738174383499
$ws.Range("M20").Value = 1.068245978258965
$ws.Range("N20").Value = 1.020564325616326

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043235483385343
$ws.Range("D21").Value = 1.043860798788945
$ws.Range("E21").Value = 1.056535473799912
$ws.Range("F21").Value = 1.06289930821326
$ws.Range("I21").Value = 1.037242751980393
$ws.Range("J21").Value = 1.049769607425781
$ws.Range("K21").Value = 1.047411925709042
$ws.Range("L21").Value = 1.060040826861133
$ws.Range("M21").Value = 1.066382132502717
$ws.Range("N21").Value = 1.020043178542782

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042053670385052
$ws.Range("D22").Value = 1.042954868146051
$ws.Range("E22").Value = 1.055334458842699
$ws.Range("F22").Value = 1.061594329375486
$ws.Range("I22").Value = 1.03693099077417
$ws.Range("J22").Value = 1.048838243346299
$ws.Range("K22").Value = 1.046637306764683
$ws.Range("L22").Value = 1.05897037016135
$ws.Range("M22").Value = 1.065207170477851
$ws.Range("N22").Value = 1.019713734008882

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042680481331241
$ws.Range("D23").Value = 1.043435408244705
$ws.Range("E23").Value = 1.055971410517318
$ws.Range("F23").Value = 1.062286380985061
$ws.Range("I23").Value = 1.037096560902049
$ws.Range("J23").Value = 1.049332309391904
$ws.Range("K23").Value = 1.047048289358135
$ws.Range("L23").Value = 1.059538153790406
$ws.Range("M23").Value = 1.065830337829937
$ws.Range("N23").Value = 1.019888548840864

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045141403950005
$ws.Range("D24").Value = 1.045320900531242
$ws.Range("E24").Value = 1.058473109437174
$ws.Range("F24").Value = 1.065005306059455
$ws.Range("I24").Value = 1.037741778086114
$ws.Range("J24").Value = 1.051270103182325
$ws.Range("K24").Value = 1.048658754223422
$ws.Range("L24").Value = 1.061766576283512
$ws.Range("M24").Value = 1.068277174645689
$ws.Range("N24").Value = 1.020573032886873

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047983698281371
$ws.Range("D25").Value = 1.047496145277952
$ws.Range("E25").Value = 1.061364473181443
$ws.Range("F25").Value = 1.068149427307477
$ws.Range("I25").Value = 1.038476980140743
$ws.Range("J25").Value = 1.053504132773835
$ws.Range("K25").Value = 1.050512364948853
$ws.Range("L25").Value = 1.064338770048592
$ws.Range("M25").Value = 1.021359738156885
